$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$val8 = @'
TC11_CDS_Filter_InstrumentModel-Illumina NovaSeq 6000_Neo4jData.xlsx
'@

$val9 = @'
TC11_CDS_Filter_InstrumentModel-Illumina NovaSeq 6000_WebData.xlsx
'@

$val10 = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina NovaSeq 6000']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p, s, collect(distinct samp.sample_id) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY `Participant ID`LIMIT 100
'@

$val11 = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina NovaSeq 6000']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@

$val12 = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina NovaSeq 6000']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,f,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN 
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER By f.file_name LIMIT 100
'@

$val13 = @'
MATCH (f:file)
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina NovaSeq 6000']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,f, s, collect(distinct samp.sample_id) as samp
RETURN
count(distinct s) AS Studies,
count(distinct p) AS Participants,
count(distinct samp) AS Samples,
count(distinct f) AS Files
'@

$val14 = @'
MATCH (f:file)
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina NovaSeq 6000']MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,f, s, collect(distinct samp.sample_id) as samp
RETURN
count(distinct s) AS Studies,
count(distinct p) AS Participants,
count(distinct samp) AS Samples,
count(distinct f) AS Files
'@


# Column D width change (78.85546875 -> 98.85546875).
# Note: this COM emulation snaps ColumnWidth to 1/6-character granularity
# (same granularity the original 78.85546875 value is subject to), so 98
# is the closest settable value to the true 98.85546875 target.
$ws.Columns.Item(4).ColumnWidth = 98

# Order of writes mirrors the target shared-string insertion order:
# filenames (D2,E2), then participant/sample/file query text (B2,B3,B4),
# then the two StatQuery variants (C2,C4); remaining cells reuse these strings.
$ws.Range("D2").Value = $val8
$ws.Range("E2").Value = $val9
$ws.Range("B2").Value = $val10
$ws.Range("B3").Value = $val11
$ws.Range("B4").Value = $val12
$ws.Range("C2").Value = $val13
$ws.Range("C4").Value = $val14

# Remaining cells reference already-introduced strings
$ws.Range("D3").Value = $val8
$ws.Range("D4").Value = $val8
$ws.Range("E3").Value = $val9
$ws.Range("E4").Value = $val9
$ws.Range("C3").Value = $val13

# Update selected cell to C4
$ws.Range("C4").Select()